# Daily Report update: append the new day's (date serial 46045 / 2026-01-23)
# depository rows to Daily_Data, and refresh the derived totals on
# Today_Summary (latest-day Eligible/Registered/Total_Stock) and
# Monthly_Stats (month-to-date Eligible/Grand_Total plus BRINK'S, INC.
# Eligible WITHDRAWN/TOTAL_TODAY) that roll up BRINK'S, INC. Eligible's
# new 699.95 withdrawal.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Daily_Data: append 22 rows (11 depositories x Registered/Eligible) for
#    the new date serial 46045.
# ---------------------------------------------------------------------------
$wsDaily = $wb.Worksheets.Item("Daily_Data")

$newDate = 46045
$startRow = 310

# Columns: Region_Type, PREV_TOTAL, RECEIVED, WITHDRAWN, NET_CHANGE, ADJUSTMENT, TOTAL_TODAY
$rows = @(
    @("ASAHI DEPOSITORY LLC Registered", 0, 0, 0, 0, 0, 0),
    @("ASAHI DEPOSITORY LLC Eligible", 0, 0, 0, 0, 0, 0),
    @("BRINK'S, INC. Registered", 87949.747, 0, 0, 0, 0, 87949.747),
    @("BRINK'S, INC. Eligible", 31278.302, 0, 699.95, -699.95, 0, 30578.352),
    @("CNT DEPOSITORY, INC. Registered", 1246.06, 0, 0, 0, 0, 1246.06),
    @("CNT DEPOSITORY, INC. Eligible", 0, 0, 0, 0, 0, 0),
    @("DELAWARE DEPOSITORY Registered", 1633.941, 0, 0, 0, 0, 1633.941),
    @("DELAWARE DEPOSITORY Eligible", 18459.584, 0, 0, 0, 0, 18459.584),
    @("HSBC BANK, USA Registered", 1394.758, 0, 0, 0, 0, 1394.758),
    @("HSBC BANK, USA Eligible", 9281.978999999999, 0, 0, 0, 0, 9281.978999999999),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 2395.448, 0, 0, 0, 0, 2395.448),
    @("INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 0, 0, 0, 0, 0, 0),
    @("JP MORGAN CHASE BANK NA Registered", 114985.579, 0, 0, 0, 0, 114985.579),
    @("JP MORGAN CHASE BANK NA Eligible", 135413.823, 0, 0, 0, 0, 135413.823),
    @("LOOMIS INTERNATIONAL (US) LLC Registered", 63745.991, 0, 0, 0, 0, 63745.991),
    @("LOOMIS INTERNATIONAL (US) LLC Eligible", 132077.206, 0, 0, 0, 0, 132077.206),
    @("MALCA-AMIT USA, LLC Registered", 395.145, 0, 0, 0, 0, 395.145),
    @("MALCA-AMIT USA, LLC Eligible", 0, 0, 0, 0, 0, 0),
    @("MANFRA, TORDELLA & BROOKES, LLC Registered", 50220.42, 0, 0, 0, 0, 50220.42),
    @("MANFRA, TORDELLA & BROOKES, LLC Eligible", 1271.373, 0, 0, 0, 0, 1271.373),
    @("STONEX PRECIOUS METALS LLC Registered", 14122.765, 0, 0, 0, 0, 14122.765),
    @("STONEX PRECIOUS METALS LLC Eligible", 16.075, 0, 0, 0, 0, 16.075)
)

$endRow = $startRow + $rows.Length - 1
$wsDaily.Range("A" + $startRow + ":A" + $endRow).EntireRow.Insert()

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $wsDaily.Cells.Item($r, 1).Value = $newDate
    $wsDaily.Cells.Item($r, 2).Value = $row[0]
    $wsDaily.Cells.Item($r, 3).Value = $row[1]
    $wsDaily.Cells.Item($r, 4).Value = $row[2]
    $wsDaily.Cells.Item($r, 5).Value = $row[3]
    $wsDaily.Cells.Item($r, 6).Value = $row[4]
    $wsDaily.Cells.Item($r, 7).Value = $row[5]
    $wsDaily.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2. Today_Summary: BRINK'S, INC. row (row 3) reflects the new Eligible
#    balance (699.95 withdrawn) and the resulting Total_Stock.
# ---------------------------------------------------------------------------
$wsToday = $wb.Worksheets.Item("Today_Summary")
$wsToday.Cells.Item(3, 2).Value = 30578.352
$wsToday.Cells.Item(3, 4).Value = 118528.099

# ---------------------------------------------------------------------------
# 3. Monthly_Stats: month-to-date (2026-01) Eligible / Grand_Total roll-up,
#    and BRINK'S, INC. Eligible's WITHDRAWN / TOTAL_TODAY detail row.
# ---------------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly_Stats")
$wsMonthly.Cells.Item(2, 2).Value = 327098.392
$wsMonthly.Cells.Item(2, 4).Value = 665188.246

$wsMonthly.Cells.Item(9, 4).Value = 1369.594
$wsMonthly.Cells.Item(9, 5).Value = 30578.352
